# Conserto do erro com o rotulo da coluna 2050 nas tabelas e
# retirada das linhas com total das tabelas.

$wb = $excel.ActiveWorkbook

# xlPasteFormats - used below to copy just the (bold/centered/bordered)
# header style from a neighbouring cell onto the fixed E1 label, so the
# corrected cell keeps looking like the other year-header cells instead of
# picking up a brand-new "Text" number-format style.
$xlPasteFormats = -4122

function Fix-YearLabel($ws, $label) {
    # E1 currently holds a stray numeric value (575.64...) where a text
    # year label belongs (matching B1/C1/D1, e.g. 2015/2030/2040). Writing
    # it with a leading apostrophe forces it to be stored as text instead
    # of being reinterpreted as a number.
    $ws.Range("E1").Value = "'" + $label
    # Re-apply D1's formatting (bold, centered, bordered header look) so
    # the cell's style matches the rest of the header row again.
    $ws.Range("D1").Copy()
    $ws.Range("E1").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# Sheets 1-3 ("Potencia Acumulada", "Geracao Periodo Medio", "Atendimento a
# Ponta") use a plain "2050" label in E1 and have a "Total" row as the last
# row (row 13) that must be removed.
$simpleLabelSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)

foreach ($name in $simpleLabelSheets) {
    $ws = $wb.Worksheets.Item($name)
    Fix-YearLabel $ws "2050"
    # Drop the trailing "Total" row.
    $ws.Rows.Item(13).Delete()
}

# Sheet 4 ("Potencia Incremental") follows the period-range label scheme
# (2015-2030, 2031-2040), so its fixed label is "2041-2050" instead of
# plain "2050". It also has a trailing "Total" row to remove.
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Fix-YearLabel $ws4 "2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet 5 ("Emissoes Totais") only needs the E1 label fix -- it has no
# Total row to remove.
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Fix-YearLabel $ws5 "2050"

# Sheet 6 ("Custo Total") has no year-label row, just a trailing "Total"
# row (row 4) to remove.
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
